# Moni.xlsx edit: rename Sheet1 -> "raw data", add a new "linear regression"
# sheet with the fitted model coefficients + a prediction row, make the new
# sheet the active tab, and move the selection on "raw data".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing sheet and change its selection.
# ---------------------------------------------------------------------
$raw = $wb.Worksheets.Item(1)
$raw.Name = "raw data"
$raw.Range("F6:L6").Select()

# ---------------------------------------------------------------------
# 2) Add the new "linear regression" sheet right after "raw data".
# ---------------------------------------------------------------------
$lr = $wb.Worksheets.Add($null, $raw)
$lr.Name = "linear regression"

# Header row: a subset of the "raw data" column headers, in model order.
# (K1's "Interstital Fluid Prediction" label is written further below,
# after the footnote, so new shared strings land at the same indices as
# the authored workbook: the footnote text first, then the label.)
$lr.Range("D1").Value = "Acetone ketones ppm"
$lr.Range("E1").Value = "Blood Pressure Top"
$lr.Range("F1").Value = "Blood Pressure Bottom"
$lr.Range("G1").Value = "Pulse"
$lr.Range("H1").Value = "Heart Rate"
$lr.Range("I1").Value = "Temperature"
$lr.Range("J1").Value = "Pulse Oxygen"

# Sample row used to demo the prediction formula.
$lr.Range("D2").Value = 0
$lr.Range("E2").Value = 125
$lr.Range("F2").Value = 77
$lr.Range("G2").Value = 106
$lr.Range("H2").Value = 100
$lr.Range("I2").Value = 97.3
$lr.Range("J2").Value = 96

# Regression coefficients (intercept + 7 slopes) from the fitted model.
$lr.Range("B3").Value = -0.0181
$lr.Range("B4").Value = -0.4819
$lr.Range("B5").Value = -14.4261
$lr.Range("B6").Value = 36.5203
$lr.Range("B7").Value = -7.5521
$lr.Range("B8").Value = 2.5154
$lr.Range("B9").Value = -14.2426
$lr.Range("B10").Value = 10.96

# Prediction formula.
$lr.Range("K2").Formula = "=B3+B4*D2+B5*E2+B6*F2+B7*G2+B8*H2+B9*I2+B10*J2"

# Footnote describing the fitted regression equation.
$lr.Range("A16").Value = "Interstitial Fluid = -0.0181 - 0.4819 * Acetone ketones ppm - 14.4261 * Blood Pressure Top + 36.5203 * Blood Pressure Bottom - 7.5521 * Pulse + 2.5154 * Heart Rate - 14.2426 * Temperature + 10.9600 * Pulse Oxygen"

# Prediction column header (written after the footnote so it becomes the
# next new shared string, matching the source sheet's string order).
$lr.Range("K1").Value = "Interstital Fluid Prediction"

# ---------------------------------------------------------------------
# 3) Formatting: bigger custom font for the footnote, green highlight
#    fill on the prediction header/value (order matters for style index
#    allocation: font style first, then fill style).
# ---------------------------------------------------------------------
$lr.Range("A16").Font.Name = "Var(--jp-code-font-family)"
$lr.Range("A16").Font.Size = 13

$lr.Range("K1:K2").Interior.Color = 5296274

$lr.Rows.Item(16).RowHeight = 17

# ---------------------------------------------------------------------
# 4) Column widths on "linear regression" to match the authored layout.
# ---------------------------------------------------------------------
$lr.Columns.Item(4).ColumnWidth = 17.666666666666668
$lr.Columns.Item(5).ColumnWidth = 15.666666666666666
$lr.Columns.Item(6).ColumnWidth = 18.833333333333332
$lr.Columns.Item(7).ColumnWidth = 4.666666666666667
$lr.Columns.Item(8).ColumnWidth = 9.166666666666666
$lr.Columns.Item(9).ColumnWidth = 10.666666666666666
$lr.Columns.Item(10).ColumnWidth = 10.833333333333332
$lr.Columns.Item(11).ColumnWidth = 21.666666666666668

# ---------------------------------------------------------------------
# 5) Make "linear regression" the active/selected tab.
# ---------------------------------------------------------------------
$lr.Range("D3").Select()
$lr.Select()
